$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.2119205298013245
$ws.Range("C2").Value = 0.4966887417218543
$ws.Range("J2").Value = 0.009933774834437087
$ws.Range("P2").Value = 0.1821192052980132
$ws.Range("S2").Value = 0.09933774834437085
$ws.Range("B3").Value = 0.0125
$ws.Range("C3").Value = 0.01875
$ws.Range("J3").Value = 0.025
$ws.Range("P3").Value = 0.725
$ws.Range("S3").Value = 0.21875
$ws.Range("J4").Value = 0.04545454545454546
$ws.Range("P4").Value = 0.5227272727272727
$ws.Range("S4").Value = 0.4318181818181818
$ws.Range("B6").Value = 0.06204379562043796
$ws.Range("D6").Value = 0.0145985401459854
$ws.Range("F6").Value = 0.05109489051094891
$ws.Range("J6").Value = 0.2883211678832117
$ws.Range("O6").Value = 0.0218978102189781
$ws.Range("Q6").Value = 0.1423357664233577
$ws.Range("R6").Value = 0.0583941605839416
$ws.Range("S6").Value = 0.3613138686131387
$ws.Range("B7").Value = 0.1650485436893204
$ws.Range("D7").Value = 0.009708737864077669
$ws.Range("F7").Value = 0.05825242718446602
$ws.Range("J7").Value = 0.1407766990291262
$ws.Range("O7").Value = 0.04854368932038835
$ws.Range("Q7").Value = 0.116504854368932
$ws.Range("R7").Value = 0.04368932038834952
$ws.Range("S7").Value = 0.4174757281553398
$ws.Range("B8").Value = 0.09070796460176991
$ws.Range("D8").Value = 0.01991150442477876
$ws.Range("F8").Value = 0.05309734513274336
$ws.Range("J8").Value = 0.1305309734513274
$ws.Range("O8").Value = 0.04424778761061947
$ws.Range("Q8").Value = 0.1349557522123894
$ws.Range("R8").Value = 0.06415929203539823
$ws.Range("S8").Value = 0.4623893805309734
$ws.Range("B9").Value = 0.0995260663507109
$ws.Range("D9").Value = 0.02369668246445497
$ws.Range("E9").Value = 0.004739336492890996
$ws.Range("F9").Value = 0.09004739336492891
$ws.Range("J9").Value = 0.1184834123222749
$ws.Range("O9").Value = 0.03791469194312796
$ws.Range("Q9").Value = 0.1421800947867299
$ws.Range("R9").Value = 0.09004739336492891
$ws.Range("S9").Value = 0.3933649289099526
$ws.Range("B10").Value = 0.1258741258741259
$ws.Range("D10").Value = 0.02272727272727273
$ws.Range("E10").Value = 0.0008741258741258741
$ws.Range("F10").Value = 0.08916083916083917
$ws.Range("J10").Value = 0.1048951048951049
$ws.Range("O10").Value = 0.03409090909090909
$ws.Range("Q10").Value = 0.1923076923076923
$ws.Range("R10").Value = 0.0472027972027972
$ws.Range("S10").Value = 0.3828671328671329
$ws.Range("G11").Value = 0.1712328767123288
$ws.Range("J11").Value = 0.08904109589041095
$ws.Range("K11").Value = 0.1917808219178082
$ws.Range("L11").Value = 0.5273972602739726
$ws.Range("S11").Value = 0.02054794520547945
$ws.Range("G12").Value = 0.6778846153846154
$ws.Range("J12").Value = 0.1490384615384615
$ws.Range("K12").Value = 0.02884615384615385
$ws.Range("L12").Value = 0.0625
$ws.Range("S12").Value = 0.08173076923076923
$ws.Range("F13").Value = 0.02040816326530612
$ws.Range("G13").Value = 0.6938775510204082
$ws.Range("J13").Value = 0.1020408163265306
$ws.Range("S13").Value = 0.1836734693877551
$ws.Range("G14").Value = 0.5
$ws.Range("J14").Value = 0.5
$ws.Range("F15").Value = 0.04705882352941176
$ws.Range("H15").Value = 0.1725490196078431
$ws.Range("I15").Value = 0.05882352941176471
$ws.Range("J15").Value = 0.2745098039215687
$ws.Range("K15").Value = 0.07058823529411765
$ws.Range("M15").Value = 0.007843137254901961
$ws.Range("N15").Value = 0.00392156862745098
$ws.Range("O15").Value = 0.05098039215686274
$ws.Range("S15").Value = 0.3137254901960784
$ws.Range("F16").Value = 0.03684210526315789
$ws.Range("H16").Value = 0.2052631578947368
$ws.Range("I16").Value = 0.08947368421052632
$ws.Range("J16").Value = 0.3421052631578947
$ws.Range("K16").Value = 0.1263157894736842
$ws.Range("M16").Value = 0.02105263157894737
$ws.Range("N16").Value = 0.005263157894736842
$ws.Range("O16").Value = 0.04210526315789474
$ws.Range("S16").Value = 0.131578947368421
$ws.Range("F17").Value = 0.02425876010781671
$ws.Range("H17").Value = 0.1482479784366577
$ws.Range("I17").Value = 0.08894878706199461
$ws.Range("J17").Value = 0.431266846361186
$ws.Range("K17").Value = 0.1051212938005391
$ws.Range("M17").Value = 0.02425876010781671
$ws.Range("O17").Value = 0.07277628032345014
$ws.Range("S17").Value = 0.1051212938005391
$ws.Range("F18").Value = 0.02362204724409449
$ws.Range("H18").Value = 0.1653543307086614
$ws.Range("I18").Value = 0.1181102362204724
$ws.Range("J18").Value = 0.3858267716535433
$ws.Range("K18").Value = 0.09448818897637795
$ws.Range("M18").Value = 0.01574803149606299
$ws.Range("O18").Value = 0.07086614173228346
$ws.Range("S18").Value = 0.1259842519685039
$ws.Range("F19").Value = 0.04338070306656694
$ws.Range("H19").Value = 0.2056843679880329
$ws.Range("I19").Value = 0.08825729244577413
$ws.Range("J19").Value = 0.337322363500374
$ws.Range("K19").Value = 0.09274495138369485
$ws.Range("M19").Value = 0.01795063575168287
$ws.Range("N19").Value = 0.002243829468960359
$ws.Range("O19").Value = 0.07180254300673149
$ws.Range("S19").Value = 0.1406133133881825
